$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellRef, $value)
    $helper = $ws.Range("ZZ1")
    $helper.NumberFormat = "@"
    $helper.Value = $value
    $helper.Copy()
    $dst = $ws.Range($cellRef)
    $dst.PasteSpecial(-4163, -4142, $false, $false)
    $helper.Clear()
}

Set-TextValue $ws "D2" '27.967.47'
Set-TextValue $ws "E2" '  +0.84%  '

Set-TextValue $ws "D3" '1.894.38'
Set-TextValue $ws "E3" '  +0.41%  '

Set-TextValue $ws "E4" '  +1.62%  '

Set-TextValue $ws "D5" '336.25'
Set-TextValue $ws "E5" '  +1.56%  '

Set-TextValue $ws "D6" '1.017'
Set-TextValue $ws "E6" '  +1.66%  '

Set-TextValue $ws "D7" '0.4696'
Set-TextValue $ws "E7" '  -0.62%  '

Set-TextValue $ws "D8" '0.3927'
Set-TextValue $ws "E8" '  -1.21%  '

Set-TextValue $ws "D9" '47.71'
Set-TextValue $ws "E9" '  -1.48%  '

Set-TextValue $ws "D10" '0.08050'
Set-TextValue $ws "E10" '  -0.17%  '

Set-TextValue $ws "D11" '1.021'
Set-TextValue $ws "E11" '  -0.74%  '

Set-TextValue $ws "D12" '21.84'
Set-TextValue $ws "E12" '  -0.21%  '

Set-TextValue $ws "D13" '1.893.85'
Set-TextValue $ws "E13" '  +0.84%  '

Set-TextValue $ws "D14" '5.971'
Set-TextValue $ws "E14" '  -0.04%  '

Set-TextValue $ws "D15" '7.119'
Set-TextValue $ws "E15" '  -1.35%  '

Set-TextValue $ws "D16" '1.020'

Set-TextValue $ws "D17" '0.06808'
Set-TextValue $ws "E17" '  +3.20%  '

Set-TextValue $ws "D18" '0.00001055'
Set-TextValue $ws "E18" '  +1.27%  '

Set-TextValue $ws "D19" '87.55'
Set-TextValue $ws "E19" '  +0.60%  '

Set-TextValue $ws "D20" '17.18'
Set-TextValue $ws "E20" '  -0.84%  '

Set-TextValue $ws "E21" '  +1.59%  '

Set-TextValue $ws "D22" '27.990.38'
Set-TextValue $ws "E22" '  +0.89%  '

Set-TextValue $ws "D23" '5.520'
Set-TextValue $ws "E23" '  -0.05%  '

Set-TextValue $ws "D24" '10.99'
Set-TextValue $ws "E24" '  -0.19%  '

Set-TextValue $ws "D25" '2.348'
Set-TextValue $ws "E25" '  +1.68%  '

Set-TextValue $ws "D26" '2.124.43'
Set-TextValue $ws "E26" '  +1.06%  '

Set-TextValue $ws "D27" '159.42'
Set-TextValue $ws "E27" '  +2.90%  '

Set-TextValue $ws "D28" '20.06'
Set-TextValue $ws "E28" '  -0.93%  '

Set-TextValue $ws "D29" '2.085'
Set-TextValue $ws "E29" '  -0.83%  '

Set-TextValue $ws "D30" '5.465'
Set-TextValue $ws "E30" '  -2.67%  '

Set-TextValue $ws "D31" '122.12'
Set-TextValue $ws "E31" '  -0.57%  '

Set-TextValue $ws "D32" '0.9711'
Set-TextValue $ws "E32" '  +0.36%  '

Set-TextValue $ws "D34" '3.675'
Set-TextValue $ws "E34" '  +1.41%  '

Set-TextValue $ws "D35" '1.400'
Set-TextValue $ws "E35" '  -5.18%  '

Set-TextValue $ws "D36" '5.372'
Set-TextValue $ws "E36" '  +1.14%  '

Set-TextValue $ws "D37" '0.06138'
Set-TextValue $ws "E37" '  +0.02%  '

Set-TextValue $ws "D38" '0.02257'
Set-TextValue $ws "E38" '  +0.07%  '

Set-TextValue $ws "D39" '1.220'
Set-TextValue $ws "E39" '  -0.28%  '

Set-TextValue $ws "D40" '8.134'
Set-TextValue $ws "E40" '  -0.66%  '

Set-TextValue $ws "D41" '0.6002'
Set-TextValue $ws "E41" '  -0.22%  '

Set-TextValue $ws "D42" '0.1891'
Set-TextValue $ws "E42" '  -0.64%  '

Set-TextValue $ws "D43" '10.34'
Set-TextValue $ws "E43" '  -0.15%  '

Set-TextValue $ws "D44" '1.271'
Set-TextValue $ws "E44" '  +0.58%  '

Set-TextValue $ws "D45" '0.5705'
Set-TextValue $ws "E45" '  +0.05%  '

Set-TextValue $ws "D46" '12.28'
Set-TextValue $ws "E46" '  -0.03%  '

Set-TextValue $ws "D47" '3.413'
Set-TextValue $ws "E47" '  +0.22%  '

Set-TextValue $ws "D48" '1.939'
Set-TextValue $ws "E48" '  -0.14%  '

Set-TextValue $ws "D49" '0.06939'
Set-TextValue $ws "E49" '  +1.60%  '

Set-TextValue $ws "D50" '114.08'
Set-TextValue $ws "E50" '  +2.96%  '

Set-TextValue $ws "D51" '1.071'
Set-TextValue $ws "E51" '  +0.08%  '
